# Update forecast values in columns C (yhat_lower) and D (yhat_upper)
# for rows 2-6, reflecting the revised model run after the TFM defense.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -4976.313244358554
$ws.Range("D2").Value = 16209.14789398284

$ws.Range("C3").Value = 636.2068445654998
$ws.Range("D3").Value = 21164.59471585847

$ws.Range("C4").Value = 11947.94006750283
$ws.Range("D4").Value = 33657.12323262377

$ws.Range("C5").Value = 2130.00157561662
$ws.Range("D5").Value = 23817.13369383052

$ws.Range("C6").Value = 3452.861499942476
$ws.Range("D6").Value = 23777.46398335794
